$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '62.611.33'
Set-TextValue $ws.Range("E2") '  -0.88%  '
Set-TextValue $ws.Range("D3") '2.571.25'
Set-TextValue $ws.Range("E3") '  -0.03%  '
Set-TextValue $ws.Range("D5") '579.82'
Set-TextValue $ws.Range("E5") '  -0.50%  '
Set-TextValue $ws.Range("D6") '143.65'
Set-TextValue $ws.Range("E6") '  -3.20%  '
Set-TextValue $ws.Range("E7") '  +0.02%  '
Set-TextValue $ws.Range("E8") '  +0.34%  '
Set-TextValue $ws.Range("E9") '  -1.34%  '
Set-TextValue $ws.Range("E10") '  -1.04%  '
Set-TextValue $ws.Range("E12") '  -1.63%  '
Set-TextValue $ws.Range("E13") '  -3.36%  '
Set-TextValue $ws.Range("D14") '3.033.42'
Set-TextValue $ws.Range("E14") '  +0.05%  '
Set-TextValue $ws.Range("D15") '62.550.90'
Set-TextValue $ws.Range("E15") '  -0.87%  '
Set-TextValue $ws.Range("E16") '  -1.25%  '
Set-TextValue $ws.Range("D17") '2.572.54'
Set-TextValue $ws.Range("E17") '  +0.16%  '
Set-TextValue $ws.Range("D18") '11.16'
Set-TextValue $ws.Range("E18") '  -2.33%  '
Set-TextValue $ws.Range("D19") '337.72'
Set-TextValue $ws.Range("E19") '  -0.97%  '
Set-TextValue $ws.Range("E20") '  -1.16%  '
Set-TextValue $ws.Range("D21") '6.64'
Set-TextValue $ws.Range("E21") '  -3.24%  '
Set-TextValue $ws.Range("E22") '  -0.04%  '
Set-TextValue $ws.Range("D23") '67.10'
Set-TextValue $ws.Range("E23") '  +1.74%  '
Set-TextValue $ws.Range("E24") '  -3.70%  '
Set-TextValue $ws.Range("E25") '  -4.16%  '
Set-TextValue $ws.Range("E26") '  +1.22%  '
Set-TextValue $ws.Range("B27") 'Aptos'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D27") '7.97'
Set-TextValue $ws.Range("E27") '  -0.79%  '
Set-TextValue $ws.Range("B28") 'Binance-PegBSC-USD'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D28") '1.00'
Set-TextValue $ws.Range("E28") '  -0.05%  '
Set-TextValue $ws.Range("E29") '  -3.63%  '
Set-TextValue $ws.Range("D30") '1.91'
Set-TextValue $ws.Range("E30") '  -2.62%  '
Set-TextValue $ws.Range("D31") '0.0₃0803'
Set-TextValue $ws.Range("E31") '  -2.98%  '
Set-TextValue $ws.Range("D32") '454.77'
Set-TextValue $ws.Range("E32") '  +3.69%  '
Set-TextValue $ws.Range("D33") '176.47'
Set-TextValue $ws.Range("E33") '  -0.57%  '
Set-TextValue $ws.Range("E34") '  +0.12%  '
Set-TextValue $ws.Range("E35") '  -0.02%  '
Set-TextValue $ws.Range("E36") '  -2.58%  '
Set-TextValue $ws.Range("D37") '18.86'
Set-TextValue $ws.Range("E37") '  -2.42%  '
Set-TextValue $ws.Range("E38") '  -2.55%  '
Set-TextValue $ws.Range("E39") '  +0.01%  '
Set-TextValue $ws.Range("D40") '1.68'
Set-TextValue $ws.Range("E40") '  -4.25%  '
Set-TextValue $ws.Range("D41") '159.15'
Set-TextValue $ws.Range("E41") '  +4.47%  '
Set-TextValue $ws.Range("E42") '  -3.81%  '
Set-TextValue $ws.Range("D43") '0.626'
Set-TextValue $ws.Range("E43") '  +3.05%  '
Set-TextValue $ws.Range("D44") '20.76'
Set-TextValue $ws.Range("E44") '  -1.93%  '
Set-TextValue $ws.Range("E45") '  -3.44%  '
Set-TextValue $ws.Range("D46") '0.0957'
Set-TextValue $ws.Range("E46") '  -1.85%  '
Set-TextValue $ws.Range("E47") '  -3.79%  '
Set-TextValue $ws.Range("D48") '17.97'
Set-TextValue $ws.Range("E48") '  -2.70%  '
Set-TextValue $ws.Range("D49") '11.41'
Set-TextValue $ws.Range("E49") '  +0.31%  '
Set-TextValue $ws.Range("E50") '  -3.96%  '
Set-TextValue $ws.Range("D51") '0.957'
Set-TextValue $ws.Range("E51") '  +3.41%  '
